$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "2972"
    B = "2025-09-10"
    C = "Erdemli"
    D = "1"
    E = "ÇAP"
    F = "AYHAN KARADAYI (K.Teknisyeni)"
}

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rng = $ws.Range("A40:F40")
    # Force the new cells to be stored as text (matching the rest of the
    # sheet, which stores every column - including numeric-looking ones -
    # as text) instead of letting Excel auto-coerce numbers/dates.
    $rng.NumberFormat = "@"
    $ws.Range("A40").Value = $newRow.A
    $ws.Range("B40").Value = $newRow.B
    $ws.Range("C40").Value = $newRow.C
    $ws.Range("D40").Value = $newRow.D
    $ws.Range("E40").Value = $newRow.E
    $ws.Range("F40").Value = $newRow.F
    # Reset back to the default style so no stray format/style index is
    # left behind on the new row (keeps it identical to the other rows).
    $rng.Style = "Normal"
}
